$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell B125: value changes from 6418.67 to 6505
$ws.Range("B125").Value = 6505

# Add new row 126 (A126 = date 05/01/2023 serial 45047, B126 = 794),
# copying the formatting (number format, font, fill, border) from row 125
# so the same style indexes are reused instead of new ones being created.
$ws.Range("A125:B125").Copy()
$ws.Range("A126:B126").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A126").Value = 45047
$ws.Range("B126").Value = 794

# Match the row height used by the other data rows.
$ws.Rows.Item(126).RowHeight = $ws.Rows.Item(125).RowHeight

# Update the sheet view's active cell/selection to B127, as in the diff.
# Use Goto with Scroll:=False so the current top-left scroll position of the
# window is disturbed as little as possible.
$excel.Goto($ws.Range("B127"), $false)
